$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-125 shift down to 11-126.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with its data (constant metadata columns
# repeated from the rest of the table, plus the new date/price figures).
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44503
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 100112037
$ws.Cells.Item(10, 7).Value = "Cebollín"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 2800
$ws.Cells.Item(10, 11).Value = 900
$ws.Cells.Item(10, 12).Value = 1000
$ws.Cells.Item(10, 13).Value = 950
$ws.Cells.Item(10, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 16).Value = 158
$ws.Cells.Item(10, 17).Value = 6
$ws.Cells.Item(10, 18).Value = "Hortaliza"
